$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; existing F (incl. its column width) shifts to G
$ws.Columns.Item(6).Insert()

# New F column width (engine snaps to 1/6 character-width granularity)
$ws.Columns.Item(6).ColumnWidth = 16.5

# D-column threshold edits
$ws.Range("D2").Value = "<24"
$ws.Range("D3").Value = "<22"
$ws.Range("D15").Value = "<14"
$ws.Range("D16").Value = "<5"

# New F (Threshold) column content
$ws.Range("F1").Value = "Threshold"
$ws.Range("F2").Value = "<28 cm"
$ws.Range("F3").Value = "<22 Pa"
$ws.Range("F4").Value = ">5cm"
$ws.Range("F5").Value = "<4000 W/sqm"
$ws.Range("F6").Value = "30-65 cm"
$ws.Range("F7").Value = "<0.8"
$ws.Range("F8").Value = "<15 cm"
$ws.Range("F9").Value = "<0.05 m/s"
$ws.Range("F10").Value = "<30 cm"
$ws.Range("F11").Value = "<0.4 m/s"
$ws.Range("F15").Value = "<18"
$ws.Range("F16").Value = "<0.5 (Pa)"
$ws.Range("F17").Value = "<5cm"
$ws.Range("F18").Value = "<4000 W/sqm"
$ws.Range("F19").Value = "35-60"
$ws.Range("F20").Value = "<0.7"
$ws.Range("F21").Value = "<15 cm"
$ws.Range("F22").Value = "<0.05 m/s"
$ws.Range("F23").Value = "<30 cm"
$ws.Range("F24").Value = "<0.4 m/s"

# Selection as in target
$ws.Range("D19:D20").Select()
